$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value2 = 8000
$ws.Range("J51").Value2 = 8000
$ws.Range("L51").Value2 = 8000
$ws.Range("N51").Value2 = -8968

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value2 = 4185.2856
$ws.Range("I62").Value2 = 3882.8333
$ws.Range("K62").Value2 = 3882.8333
$ws.Range("M62").Value2 = -3258.8333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value2 = 4185.2856
$ws.Range("I65").Value2 = 3882.8333
$ws.Range("K65").Value2 = 19414.1665
$ws.Range("M65").Value2 = -16294.1665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value2 = 6000
$ws.Range("I76").Value2 = 7000
$ws.Range("J76").Value2 = 5000
$ws.Range("K76").Value2 = 7000
$ws.Range("L76").Value2 = 5000
$ws.Range("M76").Value2 = -6685
$ws.Range("N76").Value2 = -5630

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value2 = 6000
$ws.Range("I79").Value2 = 7000
$ws.Range("J79").Value2 = 5000
$ws.Range("K79").Value2 = 7000
$ws.Range("L79").Value2 = 5000
$ws.Range("M79").Value2 = -5908
$ws.Range("N79").Value2 = -7184

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value2 = 7900
$ws.Range("I82").Value2 = 7900
$ws.Range("K82").Value2 = 23700
$ws.Range("M82").Value2 = -23294

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H85").Value2 = 7900
$ws.Range("I85").Value2 = 7900
$ws.Range("K85").Value2 = 23700
$ws.Range("M85").Value2 = -22296

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value2 = 2159.25
$ws.Range("I100").Value2 = 2210.5715
$ws.Range("J100").Value2 = 1800
$ws.Range("K100").Value2 = 2210.5715
$ws.Range("L100").Value2 = 1800
$ws.Range("M100").Value2 = -1669.5715
$ws.Range("N100").Value2 = -2882

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value2 = 4950
$ws.Range("I116").Value2 = 5000
$ws.Range("K116").Value2 = 5000
$ws.Range("M116").Value2 = -1558

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value2 = 600
$ws.Range("I118").Value2 = 600
$ws.Range("K118").Value2 = 1800
$ws.Range("M118").Value2 = -143

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value2 = 1686.75
$ws.Range("I132").Value2 = 1744.4546
$ws.Range("K132").Value2 = 5233.3638
$ws.Range("M132").Value2 = -2703.3638

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value2 = 3056.4333
$ws.Range("I137").Value2 = 2140.3635
$ws.Range("K137").Value2 = 6421.0905
$ws.Range("M137").Value2 = -3871.0905

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value2 = 4589.5713
$ws.Range("I138").Value2 = 2969.1667
$ws.Range("K138").Value2 = 8907.500100000001
$ws.Range("M138").Value2 = -3767.500100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 11611.4375
$ws.Range("I32").Value2 = 7801.091
$ws.Range("K32").Value2 = 7801.091
$ws.Range("M32").Value2 = -7514.091

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H60").Value2 = 400000
$ws.Range("I60").Value2 = 400000
$ws.Range("J60").Value2 = 0
$ws.Range("K60").Value2 = 400000
$ws.Range("L60").Value2 = 0
$ws.Range("M60").Value2 = -399267
$ws.Range("N60").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 2333
$ws.Range("I61").Value2 = 2333
$ws.Range("K61").Value2 = 2333
$ws.Range("M61").Value2 = -2121

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value2 = 3341.3333
$ws.Range("I132").Value2 = 3341.3333
$ws.Range("K132").Value2 = 10023.9999
$ws.Range("M132").Value2 = -7493.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value2 = 2333
$ws.Range("I136").Value2 = 2333
$ws.Range("K136").Value2 = 6999
$ws.Range("M136").Value2 = -4449

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value2 = 3999
$ws.Range("I134").Value2 = 3999
$ws.Range("K134").Value2 = 11997
$ws.Range("M134").Value2 = -9462

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value2 = 2648.8333
$ws.Range("I58").Value2 = 2478.6
$ws.Range("K58").Value2 = 2478.6
$ws.Range("M58").Value2 = -2275.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value2 = 1624.8572
$ws.Range("I134").Value2 = 1575.2
$ws.Range("J134").Value2 = 1749
$ws.Range("K134").Value2 = 4725.6
$ws.Range("L134").Value2 = 5247
$ws.Range("M134").Value2 = -2190.6
$ws.Range("N134").Value2 = -10317

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value2 = 2648.8333
$ws.Range("I136").Value2 = 2478.6
$ws.Range("K136").Value2 = 7435.799999999999
$ws.Range("M136").Value2 = -4885.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value2 = 0
$ws.Range("I69").Value2 = 0
$ws.Range("J69").Value2 = 0
$ws.Range("K69").Value2 = 0
$ws.Range("L69").Value2 = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value2 = 0
$ws.Range("I72").Value2 = 0
$ws.Range("J72").Value2 = 0
$ws.Range("K72").Value2 = 0
$ws.Range("L72").Value2 = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value2 = 102
$ws.Range("I92").Value2 = 102
$ws.Range("K92").Value2 = 306
$ws.Range("M92").Value2 = 942

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value2 = 25000
$ws.Range("I57").Value2 = 0
$ws.Range("J57").Value2 = 25000
$ws.Range("K57").Value2 = 0
$ws.Range("L57").Value2 = 25000
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value2 = -26640

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 3499.3333
$ws.Range("I80").Value2 = 3499
$ws.Range("J80").Value2 = 3499.5
$ws.Range("K80").Value2 = 3499
$ws.Range("L80").Value2 = 3499.5
$ws.Range("M80").Value2 = -2501
$ws.Range("N80").Value2 = -5495.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value2 = 3499.3333
$ws.Range("I83").Value2 = 3499
$ws.Range("J83").Value2 = 3499.5
$ws.Range("K83").Value2 = 17495
$ws.Range("L83").Value2 = 17497.5
$ws.Range("M83").Value2 = -12503
$ws.Range("N83").Value2 = -27481.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value2 = 7567.222
$ws.Range("I102").Value2 = 7567.222
$ws.Range("K102").Value2 = 7567.222
$ws.Range("M102").Value2 = -5945.222

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value2 = 586.25
$ws.Range("I113").Value2 = 586.25
$ws.Range("K113").Value2 = 586.25
$ws.Range("M113").Value2 = 1583.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 774.75
$ws.Range("J22").Value2 = 849.5
$ws.Range("L22").Value2 = 849.5
$ws.Range("N22").Value2 = -1439.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value2 = 774.75
$ws.Range("J27").Value2 = 849.5
$ws.Range("L27").Value2 = 849.5
$ws.Range("N27").Value2 = -1063.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value2 = 2268.375
$ws.Range("I82").Value2 = 2268.375
$ws.Range("K82").Value2 = 2268.375
$ws.Range("M82").Value2 = -1907.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value2 = 2268.375
$ws.Range("I85").Value2 = 2268.375
$ws.Range("K85").Value2 = 2268.375
$ws.Range("M85").Value2 = -1020.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value2 = 4945
$ws.Range("I107").Value2 = 4945
$ws.Range("K107").Value2 = 4945
$ws.Range("M107").Value2 = -3025

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value2 = 32428.428
$ws.Range("I132").Value2 = 34999.8
$ws.Range("J132").Value2 = 26000
$ws.Range("K132").Value2 = 104999.4
$ws.Range("L132").Value2 = 78000
$ws.Range("M132").Value2 = -102469.4
$ws.Range("N132").Value2 = -83060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value2 = 4877.5625
$ws.Range("I132").Value2 = 4217.2144
$ws.Range("J132").Value2 = 9500
$ws.Range("K132").Value2 = 12651.6432
$ws.Range("L132").Value2 = 28500
$ws.Range("M132").Value2 = -10121.6432
$ws.Range("N132").Value2 = -33560
